$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the 4 new To-Do rows (18-21)
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "Call Hiya for 16th Dec demo"
$ws.Range("C20").Value = "Not-Done"

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "Fix the TV"
$ws.Range("C21").Value = "Not-Done"

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "Purchase a white board"
$ws.Range("C22").Value = "Not-Done"

$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "Fix the white board in the demo room for 16 Dec 2023"
$ws.Range("C23").Value = "Not-Done"

# Turn on AutoFilter for the header row, which also registers the
# hidden _FilterDatabase defined name scoped to this sheet.
$ws.Range("A1:D1").AutoFilter()
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "='ToDo-List-16-Dec-2023'!`$A`$1:`$D`$1")
$fdb.Visible = $false

# Move the selection to D9, matching the saved UI state.
$ws.Range("D9").Select()
